# Apply the content edits described by the commit:
# "Themen & Ansprechpersonen leicht überarbeitet und eine Navigationsebene höher geschoben"
# (minor wording/typo fixes throughout the table, filling a few previously
#  empty "Ansprechperson"/"Website" cells with an em dash placeholder, and
#  moving the active-cell selection.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Big Team Science: fix capitalization "Erlaubt" -> "erlaubt"
$ws.Range("D2").Value = "Fördert Austausch und Zusammenarbeit, erlaubt Untersuchung von Fragestellungen, die enorme Ressourcen benötigen"

# Row 4 - Open Governance: fill previously empty Ansprechperson/Website with em dash placeholder
$ws.Range("E4").Value = "—"
$ws.Range("F4").Value = "—"

# Row 7 - Diversity: drop trailing period
$ws.Range("D7").Value = "Personengruppen werden in der Wissenschaft systematisch benachteiligt und Strukturen begünstigen Machtmissbrauch"

# Rows 9-12 - fix typo "Forschungssplanung" -> "Forschungsplanung"
$ws.Range("A9").Value = "1. Forschungsfragen und Forschungsplanung"
$ws.Range("A10").Value = "1. Forschungsfragen und Forschungsplanung"
$ws.Range("A11").Value = "1. Forschungsfragen und Forschungsplanung"
$ws.Range("A12").Value = "1. Forschungsfragen und Forschungsplanung"

# Row 10 - Sample Size Justification: remove comma before "oder"; fill Website placeholder
$ws.Range("C10").Value = "Die benötigte Menge an Beobachtungen wird im Vorhinein durch Überlegungen, Poweranalysen oder Datensimulationen festgelegt"
$ws.Range("F10").Value = "—"

# Row 11 - Theorien spezifizieren: fill Website placeholder
$ws.Range("F11").Value = "—"

# Row 15 - Open Data: "z.B." -> "z. B."; fix "msüsen" -> "müssen"
$ws.Range("C15").Value = "Forschungsdaten werden online veröffentlicht und zur Nachnutzung aufbereitet (z. B. via FAIR Kriterien und mit Codebook)"
$ws.Range("D15").Value = "Spart Kosten (z. B. weil Daten nicht erneut erhoben werden müssen), erleichtert kumulative Forschung"

# Row 16 - Anonymisierung von Daten: fill Website placeholder
$ws.Range("F16").Value = "—"

# Row 19 - Reproduktion, Replikation: clarify wording
$ws.Range("C19").Value = "Bisherige Befunde werden erneut mit denselben (Reproduktion) oder anderen Daten (Replikation) geprüft"

# Row 20 - Open Source: rewrite "vendor-lock-in" explanation; fill Website placeholder
$ws.Range("D20").Value = "Ermöglicht Qualitätskontrolle, verhindert Lock-in-Effekt (d. h. Abhängigkeit durch hohe vom Anbieter erschaffene Wechselkosten)"
$ws.Range("F20").Value = "—"

# Row 23 - Preprint Review / PCI / f1000research: rewrite Vendor-Lock-In explanation
$ws.Range("D23").Value = "Verbessert Qualitätskontrolle, ermöglicht Unabhängigkeit Forschender von kommerziellen Verlagen"

# Row 25 - Open Author Contributions: "z.B." -> "z. B."
$ws.Range("C25").Value = "Beiträge der an einer Forschungsleistung Beteiligten werden klar kommuniziert (z. B. via CRediT)"
$ws.Range("D25").Value = "Klärt Verantwortlichkeiten, erleichtert Zuordnung von Leistungen (z. B. für Forschungsevaluation), verhindert Ehrenautor*innenschaften"

# Row 26 - Open Science in der Lehre: remove comma before "oder"
$ws.Range("C26").Value = "Aspekte von Open Science werden im Rahmen von Seminaren, Vorlesungen oder entsprechenden Prüfungsleistungen diskutiert oder vorgestellt"

# Move the active cell selection up to D23 (one navigation level up)
$ws.Range("D23").Select()
